# Update the Uri values for the existing test cases (drop the leading "/")
# and record the expected status code (400) for the first test case (TCID 1,
# the one that is actually run), which represents the first successful test
# case mentioned in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ImpactData")

$ws.Range("E2").Value = "impacts/v1/impacts?startDate=2020-5-1&endDate=2022-5-1"
$ws.Range("E3").Value = "impacts/v1/impacts?startDate=2020-5-1&endDate=2022-5-1"

$ws.Range("G2").Value = "400"

$ws.Range("G2").Select()
